# ---------------------------------------------------------------------------
# "Added dec 21 classes" - adds two new worksheets (signup, myInfo) with
# sample sign-up / profile data, plus a third row on Sheet1.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Sheet1: append a third row ---------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("A3").Value = "abcde"
$sheet1.Range("B3").Value = 12345
$sheet1.Range("F27").Select()

# --- New worksheet: signup ---------------------------------------------------
$signup = $wb.Worksheets.Add($null, $sheet1)
$signup.Name = "signup"

$signup.Range("A1").Value = "firstname"
$signup.Range("B1").Value = "lastname"
$signup.Range("C1").Value = "dob"
$signup.Range("D1").Value = "emailid"

$signup.Range("A2").Value = "john"
$signup.Range("B2").Value = "peter"
$signup.Range("C2").NumberFormat = "mm-dd-yy"
$signup.Range("C2").Value = (Get-Date -Year 2022 -Month 2 -Day 22 -Hour 0 -Minute 0 -Second 0)
$signup.Range("D2").Value = "some@gmail"
$signup.Hyperlinks.Add($signup.Range("D2"), "mailto:some@gmail.com")

# propagate the date / hyperlink cell formatting down without minting new
# (duplicate) style entries - Copy + PasteSpecial(formats) reuses the xf.
$signup.Range("C2").Copy()
$signup.Range("C3:C7").PasteSpecial(-4122)
$signup.Range("D2").Copy()
$signup.Range("D3:D7").PasteSpecial(-4122)
$signup.Application.CutCopyMode = $false

for ($r = 3; $r -le 7; $r++) {
    $signup.Range("A$r").Value = "david"
    $signup.Range("B$r").Value = "r"
    $signup.Range("C$r").Value = (Get-Date -Year 2021 -Month 2 -Day 22 -Hour 0 -Minute 0 -Second 0)
    $signup.Range("D$r").Value = "s@gmail"
    $signup.Hyperlinks.Add($signup.Range("D$r"), "mailto:s@gmail.com")
}

$signup.Columns.Item(3).ColumnWidth = 9.42
$signup.Columns.Item(4).ColumnWidth = 18.6
$signup.Range("D1").EntireColumn.Select()

# --- New worksheet: myInfo ---------------------------------------------------
$myInfo = $wb.Worksheets.Add($null, $signup)
$myInfo.Name = "myInfo"

$myInfo.Range("A1").Value = "firstName"
$myInfo.Range("A1").Interior.ThemeColor = 5
$myInfo.Range("A1").Copy()
$myInfo.Range("B1:E1").PasteSpecial(-4122)
$myInfo.Application.CutCopyMode = $false

$myInfo.Range("B1").Value = "LastName"
$myInfo.Range("C1").Value = "DrivingLicenceNumber"
$myInfo.Range("D1").Value = "EID"
$myInfo.Range("E1").Value = "MaritalStatus"

foreach ($r in 2, 3) {
    $myInfo.Range("A$r").Value = "A"
    $myInfo.Range("B$r").Value = "B"
    $myInfo.Range("C$r").Value = 12345
    $myInfo.Range("D$r").Value = 9999
    $myInfo.Range("E$r").Value = "Un married"
}

$myInfo.Columns.Item(1).ColumnWidth = 15.42
$myInfo.Columns.Item(2).ColumnWidth = 12.92
$myInfo.Columns.Item(3).ColumnWidth = 25.6
$myInfo.Columns.Item(5).ColumnWidth = 17.42
$myInfo.Range("C22").Select()
